# Re-export the logical architecture matrix independent of the algorithm
# condition: the matrix now starts straight from the "Air/ Terrestrian
# Gravity" -> "Mission Mgt Subsystem" exchange, every row shifts to its
# newly computed aggregation, and the table loses its final (now
# redundant) row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Air/ Terrestrian Gravity'
$ws.Range("C2").Value = 'Mission Mgt Subsystem'
$ws.Range("D2").Value = '[''Send perceived position, attitude '', ''Send perceived position, attitude '']'
$ws.Range("E2").Value = '[''CheckWinfForce'', ''Identify Absolute Aircraft Coordinates'']'
$ws.Range("F2").Value = 'Air/ Terrestrian Gravity_to_Mission Mgt Subsystem'
$ws.Range("G2").Value = '[''wind'', ''position settings'']'

$ws.Range("B3").Value = 'Air/ Terrestrian Gravity'
$ws.Range("C3").Value = 'UAV Control Station Subsystem'
$ws.Range("D3").Value = '[''Send perceived position, attitude '']'
$ws.Range("E3").Value = '[''Monitor UAV Control'']'
$ws.Range("F3").Value = 'Air/ Terrestrian Gravity_to_UAV Control Station Subsystem'
$ws.Range("G3").Value = '[''perceived position'', ''perceived attitude'']'

$ws.Range("B4").Value = 'Aircraft'
$ws.Range("C4").Value = 'Vision Subsystem'
$ws.Range("D4").Value = '[''Send aircraft view'']'
$ws.Range("E4").Value = '[''Record photos and videos'']'
$ws.Range("F4").Value = 'Aircraft_to_Vision Subsystem'
$ws.Range("G4").Value = 'Aircraft view'

$ws.Range("B5").Value = 'Aircraft Company Database'
$ws.Range("C5").Value = 'Mission Mgt Subsystem'
$ws.Range("D5").Value = '[''Send/Receive data'']'
$ws.Range("E5").Value = '[''Build FlightPlan Relative to Aircraft Type'']'
$ws.Range("F5").Value = 'Aircraft Company Database_to_Mission Mgt Subsystem'
$ws.Range("G5").Value = 'aircraft 3D models'

$ws.Range("B6").Value = 'Airline Human Operator'
$ws.Range("C6").Value = 'Mission Mgt Subsystem'
$ws.Range("D6").Value = '[''Configurate Flight Plan'']'
$ws.Range("E6").Value = '[''Build FlightPlan Relative to Aircraft Type'']'
$ws.Range("F6").Value = 'Airline Human Operator_to_Mission Mgt Subsystem'
$ws.Range("G6").Value = 'Configuration'

$ws.Range("B7").Value = 'Flight Mgt Subsystem'
$ws.Range("C7").Value = 'Propulsion Subsystem'
$ws.Range("D7").Value = '[''Control UAV attitude'', ''Control UAV Position'']'
$ws.Range("E7").Value = '[''Generate Thrust'', ''Generate Thrust'']'
$ws.Range("F7").Value = 'Flight Mgt Subsystem_to_Propulsion Subsystem'
$ws.Range("G7").Value = '[''Attitude Request'', ''Position Request'']'

$ws.Range("B8").Value = 'Mission Mgt Subsystem'
$ws.Range("C8").Value = 'Aircraft Company Database'
$ws.Range("D8").Value = '[''Send Pictures to DB'']'
$ws.Range("E8").Value = '[''Send/Receive data'']'
$ws.Range("F8").Value = 'Mission Mgt Subsystem_to_Aircraft Company Database'
$ws.Range("G8").Value = 'Pictures'

$ws.Range("B9").Value = 'Mission Mgt Subsystem'
$ws.Range("C9").Value = 'Flight Mgt Subsystem'
$ws.Range("D9").Value = '[''Sense and Avoid Obstacles'', ''Sense and Avoid Obstacles'', ''Retrieve POI'', ''Retrieve POI'', ''Manage Mission Modes'']'
$ws.Range("E9").Value = '[''Control UAV attitude'', ''Control UAV Position'', ''Control UAV attitude'', ''Control UAV Position'', ''Emergency Landing'']'
$ws.Range("F9").Value = 'Mission Mgt Subsystem_to_Flight Mgt Subsystem'
$ws.Range("G9").Value = '[''Emergency Landing'', ''Corrected Attitude'', ''Corrected position'', ''TargetedPosition'', ''Targeted Attitude'']'

$ws.Range("B10").Value = 'Mission Mgt Subsystem'
$ws.Range("C10").Value = 'Vision Subsystem'
$ws.Range("D10").Value = '[''Manage Mission Modes'']'
$ws.Range("E10").Value = '[''Manage Photos Recording'']'
$ws.Range("F10").Value = 'Mission Mgt Subsystem_to_Vision Subsystem'
$ws.Range("G10").Value = 'Photos recording'

$ws.Range("B11").Value = 'Moving Obstacles'
$ws.Range("C11").Value = 'Mission Mgt Subsystem'
$ws.Range("D11").Value = '[''Send moving obstacle position'']'
$ws.Range("E11").Value = '[''Sense and Avoid Obstacles'']'
$ws.Range("F11").Value = 'Moving Obstacles_to_Mission Mgt Subsystem'
$ws.Range("G11").Value = 'Moving object pos'

$ws.Range("B12").Value = 'Propulsion Subsystem'
$ws.Range("C12").Value = 'UAV Pilot'
$ws.Range("D12").Value = '[''Generate Thrust'']'
$ws.Range("E12").Value = '[''Send command and position setting'']'
$ws.Range("F12").Value = 'Propulsion Subsystem_to_UAV Pilot'
$ws.Range("G12").Value = 'Total thrust'

$ws.Range("B13").Value = 'Stationary Obstacle'
$ws.Range("C13").Value = 'Mission Mgt Subsystem'
$ws.Range("D13").Value = '[''Send stationary obstacle position'']'
$ws.Range("E13").Value = '[''Sense and Avoid Obstacles'']'
$ws.Range("F13").Value = 'Stationary Obstacle_to_Mission Mgt Subsystem'
$ws.Range("G13").Value = 'stationary object position'

$ws.Range("B14").Value = 'UAV Control Station Subsystem'
$ws.Range("C14").Value = 'UAV Pilot'
$ws.Range("D14").Value = '[''Monitor UAV Control'']'
$ws.Range("E14").Value = '[''Send command and position setting'']'
$ws.Range("F14").Value = 'UAV Control Station Subsystem_to_UAV Pilot'
$ws.Range("G14").Value = 'Mission info'

$ws.Range("B15").Value = 'UAV Pilot'
$ws.Range("C15").Value = 'Mission Mgt Subsystem'
$ws.Range("D15").Value = '[''Send command and position setting'']'
$ws.Range("E15").Value = '[''Manage Mission Modes'']'
$ws.Range("F15").Value = 'UAV Pilot_to_Mission Mgt Subsystem'
$ws.Range("G15").Value = 'Start Mission'

$ws.Range("B16").Value = 'Vision Subsystem'
$ws.Range("C16").Value = 'Mission Mgt Subsystem'
$ws.Range("D16").Value = '[''Manage Photos Recording'']'
$ws.Range("E16").Value = '[''Send Pictures to DB'']'
$ws.Range("F16").Value = 'Vision Subsystem_to_Mission Mgt Subsystem'
$ws.Range("G16").Value = 'SendPictureCmd'

# Remove the now-obsolete last row (row 17) from the original 17-row table
$ws.Rows.Item(17).Delete()